$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 66.27558733333333
$ws.Range("H2").Value = 198.826762
$ws.Range("I2").Value = 0.01593739484152995
$ws.Range("J2").Value = 0.01593739484152995
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.714516333333334
$ws.Range("N2").Value = 26.143549
$ws.Range("O2").Value = 0.1832255053237971
$ws.Range("P2").Value = 0.1832255053237971
$ws.Range("Q2").Value = 577.5596883175932
$ws.Range("R2").Value = 5198.037194858338
$ws.Range("S2").Value = 0.002920137223384201
$ws.Range("T2").Value = 0.002920137223384201

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 66.27558733333333
$ws.Range("H3").Value = 198.826762
$ws.Range("I3").Value = 0.01593739484152995
$ws.Range("J3").Value = 0.01593739484152995
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 16.48752133333333
$ws.Range("N3").Value = 49.462564
$ws.Range("O3").Value = 0.3466554324170239
$ws.Range("P3").Value = 0.346655432417024
$ws.Range("Q3").Value = 1092.72016003753
$ws.Range("R3").Value = 9834.481440337768
$ws.Range("S3").Value = 0.00552478450039141
$ws.Range("T3").Value = 0.005524784500391411

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 66.27558733333333
$ws.Range("H4").Value = 198.826762
$ws.Range("I4").Value = 0.01593739484152995
$ws.Range("J4").Value = 0.01593739484152995
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.62376966666666
$ws.Range("N4").Value = 55.871309
$ws.Range("O4").Value = 0.391570739865005
$ws.Range("P4").Value = 0.391570739865005
$ws.Range("Q4").Value = 1234.301273019051
$ws.Range("R4").Value = 11108.71145717146
$ws.Range("S4").Value = 0.006240617489618595
$ws.Range("T4").Value = 0.006240617489618595

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 66.27558733333333
$ws.Range("H5").Value = 198.826762
$ws.Range("I5").Value = 0.01593739484152995
$ws.Range("J5").Value = 0.01593739484152995
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.735891666666667
$ws.Range("N5").Value = 11.207675
$ws.Range("O5").Value = 0.07854832239417409
$ws.Range("P5").Value = 0.0785483223941741
$ws.Range("Q5").Value = 247.5984144220389
$ws.Range("R5").Value = 2228.38572979835
$ws.Range("S5").Value = 0.001251855628135741
$ws.Range("T5").Value = 0.001251855628135742

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4010.868571
$ws.Range("H6").Value = 12032.605713
$ws.Range("I6").Value = 0.9644998806575645
$ws.Range("J6").Value = 0.9644998806575644
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.714516333333334
$ws.Range("N6").Value = 26.143549
$ws.Range("O6").Value = 0.1832255053237971
$ws.Range("P6").Value = 0.1832255053237971
$ws.Range("Q6").Value = 34952.77967283283
$ws.Range("R6").Value = 314575.0170554955
$ws.Range("S6").Value = 0.1767209780182242
$ws.Range("T6").Value = 0.1767209780182242

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4010.868571
$ws.Range("H7").Value = 12032.605713
$ws.Range("I7").Value = 0.9644998806575645
$ws.Range("J7").Value = 0.9644998806575644
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.48752133333333
$ws.Range("N7").Value = 49.462564
$ws.Range("O7").Value = 0.3466554324170239
$ws.Range("P7").Value = 0.346655432417024
$ws.Range("Q7").Value = 66129.28112955869
$ws.Range("R7").Value = 595163.5301660282
$ws.Range("S7").Value = 0.334349123195516
$ws.Range("T7").Value = 0.334349123195516

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4010.868571
$ws.Range("H8").Value = 12032.605713
$ws.Range("I8").Value = 0.9644998806575645
$ws.Range("J8").Value = 0.9644998806575644
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.62376966666666
$ws.Range("N8").Value = 55.871309
$ws.Range("O8").Value = 0.391570739865005
$ws.Range("P8").Value = 0.391570739865005
$ws.Range("Q8").Value = 74697.49242957648
$ws.Range("R8").Value = 672277.4318661883
$ws.Range("S8").Value = 0.3776699318687916
$ws.Range("T8").Value = 0.3776699318687915

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4010.868571
$ws.Range("H9").Value = 12032.605713
$ws.Range("I9").Value = 0.9644998806575645
$ws.Range("J9").Value = 0.9644998806575644
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.735891666666667
$ws.Range("N9").Value = 11.207675
$ws.Range("O9").Value = 0.07854832239417409
$ws.Range("P9").Value = 0.0785483223941741
$ws.Range("Q9").Value = 14984.17047049415
$ws.Range("R9").Value = 134857.5342344473
$ws.Range("S9").Value = 0.0757598475750328
$ws.Range("T9").Value = 0.07575984757503282

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.714644
$ws.Range("H10").Value = 23.143932
$ws.Range("I10").Value = 0.001855152589919057
$ws.Range("J10").Value = 0.001855152589919056
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.714516333333334
$ws.Range("N10").Value = 26.143549
$ws.Range("O10").Value = 0.1832255053237971
$ws.Range("P10").Value = 0.1832255053237971
$ws.Range("Q10").Value = 67.229391143852
$ws.Range("R10").Value = 605.064520294668
$ws.Range("S10").Value = 0.00033991127074067
$ws.Range("T10").Value = 0.00033991127074067

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 7.714644
$ws.Range("H11").Value = 23.143932
$ws.Range("I11").Value = 0.001855152589919057
$ws.Range("J11").Value = 0.001855152589919056
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 16.48752133333333
$ws.Range("N11").Value = 49.462564
$ws.Range("O11").Value = 0.3466554324170239
$ws.Range("P11").Value = 0.346655432417024
$ws.Range("Q11").Value = 127.195357529072
$ws.Range("R11").Value = 1144.758217761648
$ws.Range("S11").Value = 0.0006430987232579524
$ws.Range("T11").Value = 0.0006430987232579524

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 7.714644
$ws.Range("H12").Value = 23.143932
$ws.Range("I12").Value = 0.001855152589919057
$ws.Range("J12").Value = 0.001855152589919056
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.62376966666666
$ws.Range("N12").Value = 55.871309
$ws.Range("O12").Value = 0.391570739865005
$ws.Range("P12").Value = 0.391570739865005
$ws.Range("Q12").Value = 143.675752916332
$ws.Range("R12").Value = 1293.081776246988
$ws.Range("S12").Value = 0.0007264234721970852
$ws.Range("T12").Value = 0.0007264234721970851

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 7.714644
$ws.Range("H13").Value = 23.143932
$ws.Range("I13").Value = 0.001855152589919057
$ws.Range("J13").Value = 0.001855152589919056
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.735891666666667
$ws.Range("N13").Value = 11.207675
$ws.Range("O13").Value = 0.07854832239417409
$ws.Range("P13").Value = 0.0785483223941741
$ws.Range("Q13").Value = 28.8210742309
$ws.Range("R13").Value = 259.3896680781
$ws.Range("S13").Value = 0.0001457191237233491
$ws.Range("T13").Value = 0.0001457191237233491

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 73.63686100000001
$ws.Range("H14").Value = 220.910583
$ws.Range("I14").Value = 0.01770757191098638
$ws.Range("J14").Value = 0.01770757191098637
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.714516333333334
$ws.Range("N14").Value = 26.143549
$ws.Range("O14").Value = 0.1832255053237971
$ws.Range("P14").Value = 0.1832255053237971
$ws.Range("Q14").Value = 641.7096279198964
$ws.Range("R14").Value = 5775.386651279068
$ws.Range("S14").Value = 0.003244478811447954
$ws.Range("T14").Value = 0.003244478811447954

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 73.63686100000001
$ws.Range("H15").Value = 220.910583
$ws.Range("I15").Value = 0.01770757191098638
$ws.Range("J15").Value = 0.01770757191098637
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 16.48752133333333
$ws.Range("N15").Value = 49.462564
$ws.Range("O15").Value = 0.3466554324170239
$ws.Range("P15").Value = 0.346655432417024
$ws.Range("Q15").Value = 1214.089316657202
$ws.Range("R15").Value = 10926.80384991481
$ws.Range("S15").Value = 0.00613842599785853
$ws.Range("T15").Value = 0.006138425997858529

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 73.63686100000001
$ws.Range("H16").Value = 220.910583
$ws.Range("I16").Value = 0.01770757191098638
$ws.Range("J16").Value = 0.01770757191098637
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 18.62376966666666
$ws.Range("N16").Value = 55.871309
$ws.Range("O16").Value = 0.391570739865005
$ws.Range("P16").Value = 0.391570739865005
$ws.Range("Q16").Value = 1371.39593824035
$ws.Range("R16").Value = 12342.56344416315
$ws.Range("S16").Value = 0.006933767034397716
$ws.Range("T16").Value = 0.006933767034397714

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 73.63686100000001
$ws.Range("H17").Value = 220.910583
$ws.Range("I17").Value = 0.01770757191098638
$ws.Range("J17").Value = 0.01770757191098637
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.735891666666667
$ws.Range("N17").Value = 11.207675
$ws.Range("O17").Value = 0.07854832239417409
$ws.Range("P17").Value = 0.0785483223941741
$ws.Range("Q17").Value = 275.0993353693917
$ws.Range("R17").Value = 2475.894018324526
$ws.Range("S17").Value = 0.001390900067282179
$ws.Range("T17").Value = 0.001390900067282179
